$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Cheques")
# Re-apply the same style to A8 explicitly (no-op change) before setting value
$s = $ws.Range("A12").Style
$ws.Range("A8").Style = $s
